$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Neo4j query text to drop into A2 (also lands as a new shared string).
$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Squamous cell lung carcinoma''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

$ws.Range("A2").Value = $query

# Row grows to fit the wrapped query text.
$ws.Rows.Item(2).RowHeight = 87

# Selection moves onto the newly filled-in B2:B5 block.
$ws.Range("B2:B5").Select() | Out-Null
